$wb = $excel.ActiveWorkbook

# --- Add a new "State" column to the hotel_info sheet, between Hotel_Name and City ---
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# --- Reorder sheet tabs so "review_info" comes before "hotel_info" ---
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))
